# Atualização automática de preços de eletricidade
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 45976
$ws.Range("B2").Value = 22.99
$ws.Range("C2").Value = 19.76
$ws.Range("D2").Value = 17.6
$ws.Range("E2").Value = 16.88
$ws.Range("F2").Value = 16.5
$ws.Range("G2").Value = 17.1
$ws.Range("H2").Value = 19.64
$ws.Range("I2").Value = 21.3
$ws.Range("J2").Value = 30.73
$ws.Range("K2").Value = 28.86
$ws.Range("L2").Value = 24.32
$ws.Range("M2").Value = 18.16
$ws.Range("N2").Value = 23.04
$ws.Range("O2").Value = 24.6
$ws.Range("P2").Value = 31.76
$ws.Range("Q2").Value = 36.1
$ws.Range("R2").Value = 40.74
$ws.Range("S2").Value = 48.17
$ws.Range("T2").Value = 65.05
$ws.Range("U2").Value = 58.62
$ws.Range("V2").Value = 48.28
$ws.Range("W2").Value = 45.25
$ws.Range("X2").Value = 41.18
$ws.Range("Y2").Value = 38.52
$ws.Range("Z2").Value = 31.46
$ws.Range("AB2").Value = 53.14
$ws.Range("AD2").Value = 61.83
$ws.Range("AF2").Value = 46.76
$ws.Range("AG2").Value = "0h-13h"
